$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14, pushing existing rows 14:29 down to 15:30,
# inheriting formatting (incl. the date style on column D) from the row above.
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new weekly price record.
$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(14, 3).Value = 'Bíobío'
$ws.Cells.Item(14, 4).Value = 45128
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 'Fruta'
$ws.Cells.Item(14, 7).Value = 100104
$ws.Cells.Item(14, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(14, 9).Value = 100104003
$ws.Cells.Item(14, 10).Value = 'Membrillo'
$ws.Cells.Item(14, 11).Value = 'Champion'
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 50
$ws.Cells.Item(14, 14).Value = 12000
$ws.Cells.Item(14, 15).Value = 12000
$ws.Cells.Item(14, 16).Value = 12000
$ws.Cells.Item(14, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(14, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(14, 19).Value = 667
$ws.Cells.Item(14, 20).Value = 18
